$d = $word.ActiveDocument

$d.Content.Find.Execute("35+26=", $true, $false, $false, $false, $false, $true, 1, $false, "38+4=", 2) | Out-Null
$d.Content.Find.Execute("90-4=", $true, $false, $false, $false, $false, $true, 1, $false, "74+24=", 2) | Out-Null
$d.Content.Find.Execute("47+0=", $true, $false, $false, $false, $false, $true, 1, $false, "71+2=", 2) | Out-Null
$d.Content.Find.Execute("32+31=", $true, $false, $false, $false, $false, $true, 1, $false, "35+46=", 2) | Out-Null
$d.Content.Find.Execute("85-11=", $true, $false, $false, $false, $false, $true, 1, $false, "34+40=", 2) | Out-Null
$d.Content.Find.Execute("3+15=", $true, $false, $false, $false, $false, $true, 1, $false, "53+7=", 2) | Out-Null
$d.Content.Find.Execute("41+47=", $true, $false, $false, $false, $false, $true, 1, $false, "26-19=", 2) | Out-Null
$d.Content.Find.Execute("25+31=", $true, $false, $false, $false, $false, $true, 1, $false, "72-38=", 2) | Out-Null
$d.Content.Find.Execute("41+43=", $true, $false, $false, $false, $false, $true, 1, $false, "68+12=", 2) | Out-Null
$d.Content.Find.Execute("88-86=", $true, $false, $false, $false, $false, $true, 1, $false, "51+35=", 2) | Out-Null
$d.Content.Find.Execute("75+4=", $true, $false, $false, $false, $false, $true, 1, $false, "66-20=", 2) | Out-Null
$d.Content.Find.Execute("9+7=", $true, $false, $false, $false, $false, $true, 1, $false, "23+23=", 2) | Out-Null
$d.Content.Find.Execute("60+37=", $true, $false, $false, $false, $false, $true, 1, $false, "33-18=", 2) | Out-Null
$d.Content.Find.Execute("81+15=", $true, $false, $false, $false, $false, $true, 1, $false, "5+46=", 2) | Out-Null
$d.Content.Find.Execute("81-80=", $true, $false, $false, $false, $false, $true, 1, $false, "67-35=", 2) | Out-Null
$d.Content.Find.Execute("1+94=", $true, $false, $false, $false, $false, $true, 1, $false, "56+1=", 2) | Out-Null
$d.Content.Find.Execute("44+50=", $true, $false, $false, $false, $false, $true, 1, $false, "1+61=", 2) | Out-Null
$d.Content.Find.Execute("23+53=", $true, $false, $false, $false, $false, $true, 1, $false, "39-17=", 2) | Out-Null
$d.Content.Find.Execute("68-61=", $true, $false, $false, $false, $false, $true, 1, $false, "80-57=", 2) | Out-Null
$d.Content.Find.Execute("92-50=", $true, $false, $false, $false, $false, $true, 1, $false, "6+27=", 2) | Out-Null
$d.Content.Find.Execute("73-13=", $true, $false, $false, $false, $false, $true, 1, $false, "39+25=", 2) | Out-Null
$d.Content.Find.Execute("11+64=", $true, $false, $false, $false, $false, $true, 1, $false, "69+26=", 2) | Out-Null
$d.Content.Find.Execute("92-44=", $true, $false, $false, $false, $false, $true, 1, $false, "64-42=", 2) | Out-Null
$d.Content.Find.Execute("8+60=", $true, $false, $false, $false, $false, $true, 1, $false, "96-72=", 2) | Out-Null
$d.Content.Find.Execute("76-66=", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=", 2) | Out-Null
$d.Content.Find.Execute("51-22=", $true, $false, $false, $false, $false, $true, 1, $false, "1+30=", 2) | Out-Null
$d.Content.Find.Execute("37+12=", $true, $false, $false, $false, $false, $true, 1, $false, "66-55=", 2) | Out-Null
$d.Content.Find.Execute("83+1=", $true, $false, $false, $false, $false, $true, 1, $false, "92-5=", 2) | Out-Null
$d.Content.Find.Execute("44-26=", $true, $false, $false, $false, $false, $true, 1, $false, "51+39=", 2) | Out-Null
$d.Content.Find.Execute("56-33=", $true, $false, $false, $false, $false, $true, 1, $false, "78-65=", 2) | Out-Null
$d.Content.Find.Execute("81+3=", $true, $false, $false, $false, $false, $true, 1, $false, "87-49=", 2) | Out-Null
$d.Content.Find.Execute("72-64=", $true, $false, $false, $false, $false, $true, 1, $false, "75-18=", 2) | Out-Null
$d.Content.Find.Execute("94-71=", $true, $false, $false, $false, $false, $true, 1, $false, "40-19=", 2) | Out-Null
$d.Content.Find.Execute("9+84=", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=", 2) | Out-Null
$d.Content.Find.Execute("22+36=", $true, $false, $false, $false, $false, $true, 1, $false, "76-23=", 2) | Out-Null
$d.Content.Find.Execute("62+34=", $true, $false, $false, $false, $false, $true, 1, $false, "81-56=", 2) | Out-Null
$d.Content.Find.Execute("78-24=", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=", 2) | Out-Null
$d.Content.Find.Execute("87-11=", $true, $false, $false, $false, $false, $true, 1, $false, "43-41=", 2) | Out-Null
$d.Content.Find.Execute("64+19=", $true, $false, $false, $false, $false, $true, 1, $false, "15-2=", 2) | Out-Null
$d.Content.Find.Execute("71-16=", $true, $false, $false, $false, $false, $true, 1, $false, "39-24=", 2) | Out-Null
$d.Content.Find.Execute("76-75=", $true, $false, $false, $false, $false, $true, 1, $false, "43-20=", 2) | Out-Null
$d.Content.Find.Execute("77-39=", $true, $false, $false, $false, $false, $true, 1, $false, "76-74=", 2) | Out-Null
$d.Content.Find.Execute("47-47=", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("98-38=", $true, $false, $false, $false, $false, $true, 1, $false, "83-73=", 2) | Out-Null
$d.Content.Find.Execute("37+2=", $true, $false, $false, $false, $false, $true, 1, $false, "44-12=", 2) | Out-Null
$d.Content.Find.Execute("47+37=", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=", 2) | Out-Null
$d.Content.Find.Execute("47-24=", $true, $false, $false, $false, $false, $true, 1, $false, "38+4=", 2) | Out-Null
$d.Content.Find.Execute("5+2=", $true, $false, $false, $false, $false, $true, 1, $false, "58+34=", 2) | Out-Null
$d.Content.Find.Execute("52+39=", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=", 2) | Out-Null
$d.Content.Find.Execute("5+77=", $true, $false, $false, $false, $false, $true, 1, $false, "80-31=", 2) | Out-Null
$d.Content.Find.Execute("92-88=", $true, $false, $false, $false, $false, $true, 1, $false, "28+16=", 2) | Out-Null
$d.Content.Find.Execute("7-3=", $true, $false, $false, $false, $false, $true, 1, $false, "38+47=", 2) | Out-Null
$d.Content.Find.Execute("59+19=", $true, $false, $false, $false, $false, $true, 1, $false, "68-10=", 2) | Out-Null
$d.Content.Find.Execute("11+79=", $true, $false, $false, $false, $false, $true, 1, $false, "59-32=", 2) | Out-Null
$d.Content.Find.Execute("17+43=", $true, $false, $false, $false, $false, $true, 1, $false, "91+5=", 2) | Out-Null
$d.Content.Find.Execute("24+42=", $true, $false, $false, $false, $false, $true, 1, $false, "68+20=", 2) | Out-Null
$d.Content.Find.Execute("68-52=", $true, $false, $false, $false, $false, $true, 1, $false, "12+67=", 2) | Out-Null
$d.Content.Find.Execute("79-66=", $true, $false, $false, $false, $false, $true, 1, $false, "81+11=", 2) | Out-Null
$d.Content.Find.Execute("29-25=", $true, $false, $false, $false, $false, $true, 1, $false, "95-20=", 2) | Out-Null
$d.Content.Find.Execute("46-4=", $true, $false, $false, $false, $false, $true, 1, $false, "9+82=", 2) | Out-Null
$d.Content.Find.Execute("8+38=", $true, $false, $false, $false, $false, $true, 1, $false, "42+7=", 2) | Out-Null
$d.Content.Find.Execute("82-13=", $true, $false, $false, $false, $false, $true, 1, $false, "16+43=", 2) | Out-Null
$d.Content.Find.Execute("70-36=", $true, $false, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$d.Content.Find.Execute("5+85=", $true, $false, $false, $false, $false, $true, 1, $false, "35-2=", 2) | Out-Null
$d.Content.Find.Execute("26+32=", $true, $false, $false, $false, $false, $true, 1, $false, "0+22=", 2) | Out-Null
$d.Content.Find.Execute("36+9=", $true, $false, $false, $false, $false, $true, 1, $false, "15+0=", 2) | Out-Null
$d.Content.Find.Execute("80-13=", $true, $false, $false, $false, $false, $true, 1, $false, "25+28=", 2) | Out-Null
$d.Content.Find.Execute("5+53=", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=", 2) | Out-Null
$d.Content.Find.Execute("34-9=", $true, $false, $false, $false, $false, $true, 1, $false, "26-12=", 2) | Out-Null
$d.Content.Find.Execute("85-18=", $true, $false, $false, $false, $false, $true, 1, $false, "56-6=", 2) | Out-Null
$d.Content.Find.Execute("90-81=", $true, $false, $false, $false, $false, $true, 1, $false, "36+20=", 2) | Out-Null
$d.Content.Find.Execute("76+17=", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=", 2) | Out-Null
$d.Content.Find.Execute("30+43=", $true, $false, $false, $false, $false, $true, 1, $false, "47+10=", 2) | Out-Null
$d.Content.Find.Execute("32+32=", $true, $false, $false, $false, $false, $true, 1, $false, "78-68=", 2) | Out-Null
$d.Content.Find.Execute("19+33=", $true, $false, $false, $false, $false, $true, 1, $false, "32+56=", 2) | Out-Null
$d.Content.Find.Execute("2+29=", $true, $false, $false, $false, $false, $true, 1, $false, "43-39=", 2) | Out-Null
$d.Content.Find.Execute("95-84=", $true, $false, $false, $false, $false, $true, 1, $false, "19+26=", 2) | Out-Null
$d.Content.Find.Execute("91-53=", $true, $false, $false, $false, $false, $true, 1, $false, "42+53=", 2) | Out-Null
$d.Content.Find.Execute("88-39=", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=", 2) | Out-Null
$d.Content.Find.Execute("69-58=", $true, $false, $false, $false, $false, $true, 1, $false, "56-35=", 2) | Out-Null
$d.Content.Find.Execute("96-74=", $true, $false, $false, $false, $false, $true, 1, $false, "3+42=", 2) | Out-Null
$d.Content.Find.Execute("19+80=", $true, $false, $false, $false, $false, $true, 1, $false, "12+17=", 2) | Out-Null
$d.Content.Find.Execute("40+13=", $true, $false, $false, $false, $false, $true, 1, $false, "96-38=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $false, $false, $false, $false, $true, 1, $false, "70+14=", 2) | Out-Null
$d.Content.Find.Execute("6+92=", $true, $false, $false, $false, $false, $true, 1, $false, "43+39=", 2) | Out-Null
$d.Content.Find.Execute("99-24=", $true, $false, $false, $false, $false, $true, 1, $false, "3+51=", 2) | Out-Null
$d.Content.Find.Execute("8+43=", $true, $false, $false, $false, $false, $true, 1, $false, "7+78=", 2) | Out-Null
$d.Content.Find.Execute("64+25=", $true, $false, $false, $false, $false, $true, 1, $false, "22+69=", 2) | Out-Null
$d.Content.Find.Execute("31+22=", $true, $false, $false, $false, $false, $true, 1, $false, "24+57=", 2) | Out-Null
$d.Content.Find.Execute("76-59=", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=", 2) | Out-Null
$d.Content.Find.Execute("77+10=", $true, $false, $false, $false, $false, $true, 1, $false, "60+39=", 2) | Out-Null
$d.Content.Find.Execute("13+16=", $true, $false, $false, $false, $false, $true, 1, $false, "8+1=", 2) | Out-Null
$d.Content.Find.Execute("20+36=", $true, $false, $false, $false, $false, $true, 1, $false, "27-14=", 2) | Out-Null
$d.Content.Find.Execute("34+3=", $true, $false, $false, $false, $false, $true, 1, $false, "83-28=", 2) | Out-Null
$d.Content.Find.Execute("13+49=", $true, $false, $false, $false, $false, $true, 1, $false, "34-17=", 2) | Out-Null
$d.Content.Find.Execute("17+65=", $true, $false, $false, $false, $false, $true, 1, $false, "82-53=", 2) | Out-Null
$d.Content.Find.Execute("9+3=", $true, $false, $false, $false, $false, $true, 1, $false, "76-61=", 2) | Out-Null
$d.Content.Find.Execute("24+26=", $true, $false, $false, $false, $false, $true, 1, $false, "26+48=", 2) | Out-Null
$d.Content.Find.Execute("22-22=", $true, $false, $false, $false, $false, $true, 1, $false, "51+15=", 2) | Out-Null
$d.Content.Find.Execute("60-12=", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=", 2) | Out-Null
